$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (ALC)
$ws.Range("H40").Value = 3074.889
$ws.Range("I40").Value = 2759.8
$ws.Range("K40").Value = 2759.8
$ws.Range("M40").Value = -2584.8

# Row 41 (ALC)
$ws.Range("H41").Value = 1148.75
$ws.Range("I41").Value = 1366.6666
$ws.Range("K41").Value = 1366.6666
$ws.Range("M41").Value = -926.6666

# Row 98 (ALC)
$ws.Range("H98").Value = 38630.332
$ws.Range("I98").Value = 42697.383
$ws.Range("J98").Value = 28056
$ws.Range("K98").Value = 42697.383
$ws.Range("L98").Value = 28056
$ws.Range("M98").Value = -41199.383
$ws.Range("N98").Value = -31052

# Row 122 (ALC)
$ws.Range("H122").Value = 38630.332
$ws.Range("I122").Value = 42697.383
$ws.Range("J122").Value = 28056
$ws.Range("K122").Value = 128092.149
$ws.Range("L122").Value = 84168
$ws.Range("M122").Value = -125642.149
$ws.Range("N122").Value = -89068

# Row 132 (ALC)
$ws.Range("H132").Value = 5640.316
$ws.Range("I132").Value = 5580
$ws.Range("K132").Value = 16740
$ws.Range("M132").Value = -14210

# Row 138 (ALC)
$ws.Range("H138").Value = 4669.075
$ws.Range("I138").Value = 532.6667
$ws.Range("J138").Value = 5281.8765
$ws.Range("K138").Value = 1598.0001
$ws.Range("L138").Value = 15845.6295
$ws.Range("M138").Value = 3541.9999
$ws.Range("N138").Value = -26125.6295

$ws = $wb.Worksheets.Item("ARM")
# Row 29 (ARM)
$ws.Range("H29").Value = 9999.5
$ws.Range("J29").Value = 9999.5
$ws.Range("L29").Value = 9999.5
$ws.Range("N29").Value = -10615.5

# Row 45 (ARM)
$ws.Range("H45").Value = 8109.148
$ws.Range("I45").Value = 8685.5
$ws.Range("K45").Value = 8685.5
$ws.Range("M45").Value = -8308.5

# Row 61 (ARM)
$ws.Range("H61").Value = 8526.666999999999
$ws.Range("I61").Value = 13580.5
$ws.Range("J61").Value = 4483.6
$ws.Range("K61").Value = 13580.5
$ws.Range("L61").Value = 4483.6
$ws.Range("M61").Value = -13368.5
$ws.Range("N61").Value = -4907.6

# Row 110 (ARM)
$ws.Range("H110").Value = 1818.8572
$ws.Range("I110").Value = 1818.8572
$ws.Range("K110").Value = 1818.8572
$ws.Range("M110").Value = 226.1428000000001

# Row 122 (ARM)
$ws.Range("H122").Value = 1336766.8
$ws.Range("I122").Value = 2855.72
$ws.Range("J122").Value = 5505239
$ws.Range("K122").Value = 8567.16
$ws.Range("L122").Value = 16515717
$ws.Range("M122").Value = -6117.16
$ws.Range("N122").Value = -16520617

# Row 133 (ARM)
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# Row 136 (ARM)
$ws.Range("H136").Value = 8526.666999999999
$ws.Range("I136").Value = 13580.5
$ws.Range("J136").Value = 4483.6
$ws.Range("K136").Value = 40741.5
$ws.Range("L136").Value = 13450.8
$ws.Range("M136").Value = -38191.5
$ws.Range("N136").Value = -18550.8

$ws = $wb.Worksheets.Item("BSM")
# Row 12 (BSM)
$ws.Range("H12").Value = 6377.75
$ws.Range("I12").Value = 5505
$ws.Range("J12").Value = 6668.6665
$ws.Range("K12").Value = 5505
$ws.Range("L12").Value = 6668.6665
$ws.Range("M12").Value = -5337
$ws.Range("N12").Value = -7004.6665

# Row 134 (BSM)
$ws.Range("H134").Value = 2987
$ws.Range("I134").Value = 2987
$ws.Range("K134").Value = 8961
$ws.Range("M134").Value = -6426

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (CRP)
$ws.Range("H22").Value = 722.2069
$ws.Range("J22").Value = 742.3333
$ws.Range("L22").Value = 742.3333
$ws.Range("N22").Value = -1442.3333

# Row 31 (CRP)
$ws.Range("H31").Value = 2627.5715
$ws.Range("I31").Value = 1637.4839
$ws.Range("K31").Value = 1637.4839
$ws.Range("M31").Value = -1342.4839

# Row 34 (CRP)
$ws.Range("H34").Value = 2627.5715
$ws.Range("I34").Value = 1637.4839
$ws.Range("K34").Value = 1637.4839
$ws.Range("M34").Value = -1435.4839

# Row 99 (CRP)
$ws.Range("H99").Value = 514863.7
$ws.Range("I99").Value = 855537.3
$ws.Range("K99").Value = 855537.3
$ws.Range("M99").Value = -854039.3

# Row 126 (CRP)
$ws.Range("H126").Value = 514863.7
$ws.Range("I126").Value = 855537.3
$ws.Range("K126").Value = 2566611.9
$ws.Range("M126").Value = -2564141.9

# Row 132 (CRP)
$ws.Range("H132").Value = 32900.332
$ws.Range("J132").Value = 105403.5
$ws.Range("L132").Value = 316210.5
$ws.Range("N132").Value = -321270.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 386331.22
$ws.Range("I5").Value = 1767.4
$ws.Range("J5").Value = 910736.4399999999
$ws.Range("K5").Value = 5302.200000000001
$ws.Range("L5").Value = 2732209.32
$ws.Range("M5").Value = -5190.200000000001
$ws.Range("N5").Value = -2732433.32

# Row 12 (CUL)
$ws.Range("H12").Value = 122.6
$ws.Range("I12").Value = 112
$ws.Range("K12").Value = 336
$ws.Range("M12").Value = -163

# Row 86 (CUL)
$ws.Range("H86").Value = 1447.375
$ws.Range("I86").Value = 1600.3334
$ws.Range("J86").Value = 1355.6
$ws.Range("K86").Value = 4801.0002
$ws.Range("L86").Value = 4066.8
$ws.Range("M86").Value = -3615.0002
$ws.Range("N86").Value = -6438.799999999999

# Row 89 (CUL)
$ws.Range("H89").Value = 1447.375
$ws.Range("I89").Value = 1600.3334
$ws.Range("J89").Value = 1355.6
$ws.Range("K89").Value = 14403.0006
$ws.Range("L89").Value = 12200.4
$ws.Range("M89").Value = -8475.000599999999
$ws.Range("N89").Value = -24056.4

# Row 107 (CUL)
$ws.Range("H107").Value = 518.63635
$ws.Range("I107").Value = 203
$ws.Range("J107").Value = 550.2
$ws.Range("K107").Value = 609
$ws.Range("L107").Value = 1650.6
$ws.Range("M107").Value = 1311
$ws.Range("N107").Value = -5490.6

# Row 133 (CUL)
$ws.Range("H133").Value = 14000
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

# Row 135 (CUL)
$ws.Range("H135").Value = 386331.22
$ws.Range("I135").Value = 1767.4
$ws.Range("J135").Value = 910736.4399999999
$ws.Range("K135").Value = 15906.6
$ws.Range("L135").Value = 8196627.959999999
$ws.Range("M135").Value = -13371.6
$ws.Range("N135").Value = -8201697.959999999

# Row 140 (CUL)
$ws.Range("H140").Value = 4455.636
$ws.Range("I140").Value = 4455.636
$ws.Range("K140").Value = 13366.908
$ws.Range("M140").Value = -8186.908000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (GSM)
$ws.Range("H122").Value = 31333.455
$ws.Range("I122").Value = 33682.57
$ws.Range("K122").Value = 101047.71
$ws.Range("M122").Value = -98597.70999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 28713.475
$ws.Range("I7").Value = 32447.5
$ws.Range("K7").Value = 32447.5
$ws.Range("M7").Value = -32335.5

# Row 40 (LTW)
$ws.Range("H40").Value = 30986.646
$ws.Range("I40").Value = 34829.54
$ws.Range("K40").Value = 34829.54
$ws.Range("M40").Value = -34693.54

# Row 82 (LTW)
$ws.Range("H82").Value = 4249
$ws.Range("J82").Value = 2600
$ws.Range("L82").Value = 2600
$ws.Range("N82").Value = -3322

# Row 85 (LTW)
$ws.Range("H85").Value = 4249
$ws.Range("J85").Value = 2600
$ws.Range("L85").Value = 2600
$ws.Range("N85").Value = -5096

# Row 122 (LTW)
$ws.Range("H122").Value = 4456.5
$ws.Range("I122").Value = 3499.25
$ws.Range("J122").Value = 4839.4
$ws.Range("K122").Value = 10497.75
$ws.Range("L122").Value = 14518.2
$ws.Range("M122").Value = -8047.75
$ws.Range("N122").Value = -19418.2

# Row 126 (LTW)
$ws.Range("H126").Value = 28713.475
$ws.Range("I126").Value = 32447.5
$ws.Range("K126").Value = 97342.5
$ws.Range("M126").Value = -94872.5

# Row 132 (LTW)
$ws.Range("H132").Value = 1574077.5
$ws.Range("I132").Value = 2202020.8
$ws.Range("K132").Value = 6606062.399999999
$ws.Range("M132").Value = -6603532.399999999

# Row 136 (LTW)
$ws.Range("H136").Value = 9407.5
$ws.Range("I136").Value = 4280
$ws.Range("J136").Value = 13070
$ws.Range("K136").Value = 12840
$ws.Range("L136").Value = 39210
$ws.Range("M136").Value = -10290
$ws.Range("N136").Value = -44310

$ws = $wb.Worksheets.Item("WVR")
# Row 29 (WVR)
$ws.Range("H29").Value = 31717.5
$ws.Range("I29").Value = 37623.332
$ws.Range("K29").Value = 37623.332
$ws.Range("M29").Value = -37333.332

# Row 96 (WVR)
$ws.Range("H96").Value = 1774.25
$ws.Range("I96").Value = 1308.2
$ws.Range("J96").Value = 2107.1428
$ws.Range("K96").Value = 1308.2
$ws.Range("L96").Value = 2107.1428
$ws.Range("M96").Value = 64.79999999999995
$ws.Range("N96").Value = -4853.1428

# Row 100 (WVR)
$ws.Range("H100").Value = 28128.363
$ws.Range("I100").Value = 16254
$ws.Range("J100").Value = 103332.664
$ws.Range("K100").Value = 32508
$ws.Range("L100").Value = 206665.328
$ws.Range("M100").Value = -31967
$ws.Range("N100").Value = -207747.328

# Row 122 (WVR)
$ws.Range("H122").Value = 24807.64
$ws.Range("I122").Value = 2459.55
$ws.Range("K122").Value = 7378.650000000001
$ws.Range("M122").Value = -4928.650000000001

# Row 124 (WVR)
$ws.Range("H124").Value = 77775
$ws.Range("J124").Value = 77775
$ws.Range("L124").Value = 77775
$ws.Range("N124").Value = -87595

# Row 126 (WVR)
$ws.Range("H126").Value = 22342.682
$ws.Range("I126").Value = 29310
$ws.Range("J126").Value = 7412.7144
$ws.Range("K126").Value = 87930
$ws.Range("L126").Value = 22238.1432
$ws.Range("M126").Value = -85460
$ws.Range("N126").Value = -27178.1432
